# 9th Stab - Cosmetic Changes
# Insert two new weekly-rank columns ("Jun_15" and "Jun_17") right after the
# company-name column, pushing the existing "Jun_13" / "Jun_10" columns two
# slots to the right. The new columns are filled with the same default "UN"
# marker used throughout the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns B (Jun_13) and C (Jun_10) two places to the right,
# opening up fresh columns B and C.
$ws.Columns("B:C").Insert()

# Populate the new header cells (row 1). Set C1 before B1 so the shared
# string table records "Jun_15" ahead of "Jun_17", matching insertion order.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the new columns' data rows with the same default "UN" value used by
# every other as-yet-unranked week.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Re-apply the column width so the two freshly-inserted columns (and the
# column that was pushed into the old width-8 slot) keep a consistent,
# explicit width like the rest of the sheet.
$ws.Columns("C:E").ColumnWidth = 7.1
